$d = $word.ActiveDocument

# Update the date/weekday heading at the top of the document.
$d.Content.Find.Execute("2025-05-19 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-20 Tuesday", 2) | Out-Null

# Update each arithmetic expression cell in the table by position, since
# several "old" cell values are not unique within the document (e.g. "13-7="
# appears twice), so a global Find/Replace would be ambiguous.
$replacements = @(
    @{ R = 1; C = 1; New = "4+80=" },
    @{ R = 1; C = 2; New = "72+4=" },
    @{ R = 1; C = 3; New = "87-83=" },
    @{ R = 1; C = 4; New = "60+12=" },
    @{ R = 1; C = 5; New = "20-6=" },
    @{ R = 2; C = 1; New = "13+46=" },
    @{ R = 2; C = 2; New = "24+27=" },
    @{ R = 2; C = 3; New = "60-43=" },
    @{ R = 2; C = 4; New = "14+60=" },
    @{ R = 2; C = 5; New = "87+5=" },
    @{ R = 3; C = 1; New = "0+74=" },
    @{ R = 3; C = 2; New = "1+76=" },
    @{ R = 3; C = 3; New = "9+2=" },
    @{ R = 3; C = 4; New = "32-3=" },
    @{ R = 3; C = 5; New = "38+18=" },
    @{ R = 4; C = 1; New = "54+4=" },
    @{ R = 4; C = 2; New = "75-28=" },
    @{ R = 4; C = 3; New = "81-30=" },
    @{ R = 4; C = 4; New = "16-15=" },
    @{ R = 4; C = 5; New = "54-18=" },
    @{ R = 5; C = 1; New = "64-38=" },
    @{ R = 5; C = 2; New = "23+15=" },
    @{ R = 5; C = 3; New = "83+9=" },
    @{ R = 5; C = 4; New = "81+15=" },
    @{ R = 5; C = 5; New = "10+1=" },
    @{ R = 6; C = 1; New = "16-5=" },
    @{ R = 6; C = 2; New = "71-65=" },
    @{ R = 6; C = 3; New = "99-73=" },
    @{ R = 6; C = 4; New = "73-15=" },
    @{ R = 6; C = 5; New = "38+22=" },
    @{ R = 7; C = 1; New = "20+5=" },
    @{ R = 7; C = 2; New = "81-56=" },
    @{ R = 7; C = 3; New = "55+27=" },
    @{ R = 7; C = 4; New = "48+25=" },
    @{ R = 7; C = 5; New = "81-76=" },
    @{ R = 8; C = 1; New = "47-25=" },
    @{ R = 8; C = 2; New = "90-15=" },
    @{ R = 8; C = 3; New = "66-66=" },
    @{ R = 8; C = 4; New = "32+60=" },
    @{ R = 8; C = 5; New = "27+4=" },
    @{ R = 9; C = 1; New = "93-62=" },
    @{ R = 9; C = 2; New = "15+80=" },
    @{ R = 9; C = 3; New = "2+94=" },
    @{ R = 9; C = 4; New = "30-8=" },
    @{ R = 9; C = 5; New = "90-51=" },
    @{ R = 10; C = 1; New = "98-61=" },
    @{ R = 10; C = 2; New = "39+49=" },
    @{ R = 10; C = 3; New = "51-3=" },
    @{ R = 10; C = 4; New = "89+5=" },
    @{ R = 10; C = 5; New = "46+43=" },
    @{ R = 11; C = 1; New = "39-17=" },
    @{ R = 11; C = 2; New = "17+25=" },
    @{ R = 11; C = 3; New = "4+10=" },
    @{ R = 11; C = 4; New = "70-22=" },
    @{ R = 11; C = 5; New = "85-77=" },
    @{ R = 12; C = 1; New = "77+20=" },
    @{ R = 12; C = 2; New = "40+35=" },
    @{ R = 12; C = 3; New = "73-5=" },
    @{ R = 12; C = 4; New = "90+1=" },
    @{ R = 12; C = 5; New = "80+0=" },
    @{ R = 13; C = 1; New = "43-20=" },
    @{ R = 13; C = 2; New = "93+3=" },
    @{ R = 13; C = 3; New = "56+0=" },
    @{ R = 13; C = 4; New = "0+32=" },
    @{ R = 13; C = 5; New = "6+86=" },
    @{ R = 14; C = 1; New = "29+46=" },
    @{ R = 14; C = 2; New = "3+0=" },
    @{ R = 14; C = 3; New = "94-50=" },
    @{ R = 14; C = 4; New = "19+74=" },
    @{ R = 14; C = 5; New = "86-56=" },
    @{ R = 15; C = 1; New = "27+13=" },
    @{ R = 15; C = 2; New = "41+16=" },
    @{ R = 15; C = 3; New = "1+18=" },
    @{ R = 15; C = 4; New = "77-32=" },
    @{ R = 15; C = 5; New = "16+82=" },
    @{ R = 16; C = 1; New = "53-47=" },
    @{ R = 16; C = 2; New = "9+82=" },
    @{ R = 16; C = 3; New = "31-28=" },
    @{ R = 16; C = 4; New = "25+68=" },
    @{ R = 16; C = 5; New = "41+48=" },
    @{ R = 17; C = 1; New = "28+34=" },
    @{ R = 17; C = 2; New = "92-56=" },
    @{ R = 17; C = 3; New = "63-8=" },
    @{ R = 17; C = 4; New = "92-30=" },
    @{ R = 17; C = 5; New = "88-83=" },
    @{ R = 18; C = 1; New = "13+86=" },
    @{ R = 18; C = 2; New = "54-23=" },
    @{ R = 18; C = 3; New = "13+10=" },
    @{ R = 18; C = 4; New = "52-9=" },
    @{ R = 18; C = 5; New = "82-31=" },
    @{ R = 19; C = 1; New = "0+28=" },
    @{ R = 19; C = 2; New = "82-1=" },
    @{ R = 19; C = 3; New = "2+44=" },
    @{ R = 19; C = 4; New = "46-43=" },
    @{ R = 19; C = 5; New = "61-55=" },
    @{ R = 20; C = 1; New = "30+28=" },
    @{ R = 20; C = 2; New = "67+10=" },
    @{ R = 20; C = 3; New = "33+12=" },
    @{ R = 20; C = 4; New = "28+54=" },
    @{ R = 20; C = 5; New = "65+12=" }
)

$table = $d.Tables(1)
foreach ($item in $replacements) {
    $cell = $table.Cell($item.R, $item.C)
    $cell.Range.Text = $item.New
}
